$wb = $excel.ActiveWorkbook

# --- Populate row 2 of the whed_inst sheet ---
$wsInst = $wb.Worksheets.Item("whed_inst")
$wsInst.Range("A2").Value = 62
$wsInst.Range("B2").Value = "IAU-000062"
$wsInst.Range("C2").Value = "Acorn Institution"
$wsInst.Range("D2").Value = "Acorn Institute"
$wsInst.Range("E2").Value = "Oak Academy"

# Set selection on whed_inst to E2 and make it the active/selected tab
$wsInst.Select()
$wsInst.Range("E2").Select()

# --- Update selection on ext_inst sheet ---
$wsExt = $wb.Worksheets.Item("ext_inst")
$wsExt.Select()
$wsExt.Range("B2:D2").Select()

# Re-activate whed_inst as the final active sheet (matches activeTab=1 in workbook.xml)
$wsInst.Select()
